# The original report generator apparently built each "line item" row by
# concatenating strings into a shared buffer (string += ...) instead of
# writing each field with buffer.write(...). That bug caused the detail
# rows (item code / invoice no. / subtotal / item name) to be written one
# column too far to the left, clobbering the "C" column. The fix shifts
# those four values (C:F) one column to the right (D:G) on every affected
# detail row, leaving C blank (but still carrying its original style).
#
# xlPasteFormats constant used below to carry a source cell's number
# format / font / fill / border along when a value moves into a cell that
# previously held a different style.
$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-RightCDEFG($row) {
    $cC = $ws.Cells.Item($row, 3)
    $cD = $ws.Cells.Item($row, 4)
    $cE = $ws.Cells.Item($row, 5)
    $cF = $ws.Cells.Item($row, 6)
    $cG = $ws.Cells.Item($row, 7)

    # Snapshot the current values before anything is overwritten.
    $vC = $cC.Value()
    $vD = $cD.Value()
    $vE = $cE.Value()
    $vF = $cF.Value()

    # Work from right to left so we never clobber a value before it has
    # been read into its new home. Each destination first inherits the
    # source cell's formatting, then gets its value.
    $cF.Copy()
    $cG.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $cG.Value = $vF

    $cE.Copy()
    $cF.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $cF.Value = $vE

    $cD.Copy()
    $cE.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $cE.Value = $vD

    $cC.Copy()
    $cD.PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0
    $cD.Value = $vC

    # C keeps its original style but is now blank.
    $cC.Value = ""
}

$rows = @(4, 5, 9, 13, 17, 21, 25, 29, 30, 31, 32, 36, 37, 38, 42, 43, 47, 51)
foreach ($r in $rows) {
    Shift-RightCDEFG $r
}
